$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage so numeric-looking
# strings (e.g. "35.60") keep their original text formatting/precision
# and are not silently coerced to the Number type by Excel.
function Set-TextCellValue($range, [string]$text) {
    $origStyle = $range.Style
    $range.NumberFormat = '@'
    $range.Value = $text
    $range.Style = $origStyle
}

# Cell value updates per diff (coin price/volume refresh)
$ws.Range('D2').Value = '42.909.84'
$ws.Range('D3').Value = '2.301.45'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextCellValue $ws.Range('D5') '299.81'
$ws.Range('E5').Value = '  -0.82%  '
Set-TextCellValue $ws.Range('D6') '97.12'
$ws.Range('E6').Value = '  -1.79%  '
$ws.Range('E7').Value = '  -1.60%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -3.21%  '
Set-TextCellValue $ws.Range('D10') '35.60'
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('E12').Value = '  +1.02%  '
Set-TextCellValue $ws.Range('D13') '17.89'
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('E14').Value = '  -2.55%  '
$ws.Range('D15').Value = '2.657.07'
$ws.Range('E15').Value = '  -0.57%  '
$ws.Range('D16').Value = '2.308.33'
$ws.Range('E16').Value = '  +2.46%  '
$ws.Range('D18').Value = '42.834.39'
$ws.Range('E18').Value = '  -0.35%  '
Set-TextCellValue $ws.Range('D19') '12.79'
$ws.Range('E19').Value = '  -6.00%  '
$ws.Range('E20').Value = '  -0.93%  '
Set-TextCellValue $ws.Range('D21') '6.03'
$ws.Range('E21').Value = '  -2.51%  '
$ws.Range('E22').Value = '  -0.52%  '
Set-TextCellValue $ws.Range('D23') '240.10'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E24').Value = '  -1.57%  '
$ws.Range('E25').Value = '  +0.11%  '
Set-TextCellValue $ws.Range('D26') '2.42'
$ws.Range('E26').Value = '  -1.17%  '
$ws.Range('E27').Value = '  +0.06%  '
Set-TextCellValue $ws.Range('D28') '25.46'
$ws.Range('E28').Value = '  +2.11%  '
Set-TextCellValue $ws.Range('D29') '165.33'
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('E30').Value = '  -1.32%  '
Set-TextCellValue $ws.Range('D31') '9.04'
$ws.Range('E31').Value = '  -1.58%  '
$ws.Range('E32').Value = '  -1.43%  '
Set-TextCellValue $ws.Range('D33') '4.94'
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('E35').Value = '  -3.79%  '
Set-TextCellValue $ws.Range('D36') '16.93'
$ws.Range('E36').Value = '  -8.04%  '
$ws.Range('E37').Value = '  -1.28%  '
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('E39').Value = '  -1.40%  '
Set-TextCellValue $ws.Range('D40') '1.75'
$ws.Range('E40').Value = '  -2.98%  '
$ws.Range('E41').Value = '  -1.37%  '
$ws.Range('E42').Value = '  -2.01%  '
$ws.Range('D43').Value = '2.013.15'
$ws.Range('E43').Value = '  +0.79%  '
Set-TextCellValue $ws.Range('D44') '0.0281'
$ws.Range('E45').Value = '  +0.39%  '
$ws.Range('E46').Value = '  -1.14%  '
Set-TextCellValue $ws.Range('D47') '17.39'
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('E48').Value = '  -1.85%  '
Set-TextCellValue $ws.Range('D49') '53.49'
$ws.Range('E49').Value = '  -2.29%  '
$ws.Range('D50').Value = '2.525.58'
$ws.Range('E50').Value = '  -0.56%  '
Set-TextCellValue $ws.Range('D51') '72.06'
$ws.Range('E51').Value = '  -2.54%  '
